# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Price (col D) and Volume(1h) (col E) are stored as plain text in the sheet,
# so numeric-looking prices are written with a leading "'" to force Excel to
# keep them as text instead of silently re-typing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.785.26"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "'3.464.99"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'577.62"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'147.81"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("D7").Value = "'3.465.70"
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("D10").Value = "'7.68"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "'0.399"
$ws.Range("E12").Value = "  +3.84%  "
$ws.Range("D13").Value = "'4.055.28"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "'29.73"
$ws.Range("E14").Value = "  +6.37%  "
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "'3.464.00"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").Value = "'0.0000171"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "'62.801.49"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").Value = "'14.27"
$ws.Range("E20").Value = "  +4.95%  "
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "'388.79"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'0.556"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").Value = "'74.54"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'3.603.05"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("E28").Value = "  -5.82%  "
$ws.Range("D29").Value = "'7.53"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").Value = "'2.13"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").Value = "'23.57"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("E37").Value = "  +3.59%  "
$ws.Range("D38").Value = "'31.66"
$ws.Range("E38").Value = "  +18.82%  "
$ws.Range("D41").Value = "'3.501.26"
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("D42").Value = "'0.0753"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").Value = "'0.799"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("D44").Value = "'42.34"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  +2.94%  "
$ws.Range("E47").Value = "  +3.83%  "
$ws.Range("D48").Value = "'2.596.93"
$ws.Range("E48").Value = "  +5.54%  "
$ws.Range("E49").Value = "  +11.69%  "
$ws.Range("D50").Value = "'22.82"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "'6.70"
$ws.Range("E51").Value = "  +0.80%  "

# Rows 39/40 swapped rank: Monero moved up, ImmutableX moved down.
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'170.34"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'1.57"
$ws.Range("E40").Value = "  +6.56%  "
